# Update Release-Notes.xlsx - Folder inventory updated on Fri Jun 13 07:37:07 UTC 2025

$wb = $excel.ActiveWorkbook

# --- Sheet: Folder Inventory ---
$wsFolderInventory = $wb.Worksheets.Item("Folder Inventory")
$wsFolderInventory.Range("C2").Value = "2025-06-13 13:06:50 +0530"

# --- Sheet: Metadata ---
$wsMetadata = $wb.Worksheets.Item("Metadata")
$wsMetadata.Range("B3").Value = "2025-06-13 07:37:07 UTC"
# B5 ("Workflow Run") holds a numeric-looking value but must stay a text
# cell (as in the source file). Force text via NumberFormat, set the value,
# then clear the format again so no stray style is left behind on the cell.
$wsMetadata.Range("B5").NumberFormat = "@"
$wsMetadata.Range("B5").Value = "3"
$wsMetadata.Range("B5").ClearFormats()

# --- Sheet: Summary ---
$wsSummary = $wb.Worksheets.Item("Summary")
$wsSummary.Range("B5").Value = "2025-06-13 13:06:50 +0530"
